# PR_TestPlan.xlsx - "adding Method on documents"
# Updates the Test Plan sheet: test results marked Passed, testers assigned
# (Method/Assigned-to column), text fixes, and image counts 15 -> 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")
$wsDefects = $wb.Worksheets.Item("Defects")

# --- Automotive section (rows 3-7) ---
$ws.Range("F3").Value = "Passed"
$ws.Range("G3").Value = "Angelica"
$ws.Range("D4").Value = "User must be able to upload photo under this section"
$ws.Range("F4").Value = "Passed"
$ws.Range("G4").Value = "Angelica"
$ws.Range("E5").Value = "10 test images"
$ws.Range("F5").Value = "Passed"
$ws.Range("G5").Value = "Angelica"
$ws.Range("E6").Value = "10 new images"
$ws.Range("F6").Value = "Passed"
$ws.Range("G6").Value = "Angelica"
$ws.Range("F7").Value = "Passed"
$ws.Range("G7").Value = "Angelica"

# --- Tools & Hardware section (rows 8-12) ---
$ws.Range("F8").Value = "Passed"
$ws.Range("G8").Value = "Jonatas"
$ws.Range("D9").Value = "User must be able to upload photo under this section"
$ws.Range("F9").Value = "Passed"
$ws.Range("G9").Value = "Jonatas"
$ws.Range("E10").Value = "10 test images"
$ws.Range("F10").Value = "Passed"
$ws.Range("G10").Value = "Jonatas"
$ws.Range("E11").Value = "10 new images"
$ws.Range("F11").Value = "Passed"
$ws.Range("G11").Value = "Jonatas"
$ws.Range("F12").Value = "Passed"
$ws.Range("G12").Value = "Jonatas"

# --- Home & Pets section (rows 13-17) ---
$ws.Range("F13").Value = "Passed"
$ws.Range("G13").Value = "Ajay"
$ws.Range("D14").Value = "User must be able to upload photo under this section"
$ws.Range("F14").Value = "Passed"
$ws.Range("G14").Value = "Ajay"
$ws.Range("E15").Value = "10 test images"
$ws.Range("F15").Value = "Passed"
$ws.Range("G15").Value = "Ajay"
$ws.Range("E16").Value = "10 new images"
$ws.Range("F16").Value = "Passed"
$ws.Range("G16").Value = "Ajay"
$ws.Range("F17").Value = "Passed"
$ws.Range("G17").Value = "Ajay"

# --- Sports & Recreation section (rows 18-22) ---
$ws.Range("F18").Value = "Passed"
$ws.Range("G18").Value = "Akshita"
$ws.Range("D19").Value = "User must be able to upload photo under this section"
$ws.Range("F19").Value = "Passed"
$ws.Range("G19").Value = "Akshita"
$ws.Range("E20").Value = "10 test images"
$ws.Range("F20").Value = "Passed"
$ws.Range("G20").Value = "Akshita"
$ws.Range("E21").Value = "10 new images"
$ws.Range("F21").Value = "Passed"
$ws.Range("G21").Value = "Akshita"
$ws.Range("F22").Value = "Passed"
$ws.Range("G22").Value = "Akshita"

# --- Outdoor Living section (rows 23-27) ---
$ws.Range("F23").Value = "Passed"
$ws.Range("G23").Value = "Jonatas"
$ws.Range("D24").Value = "User must be able to upload photo under this section"
$ws.Range("F24").Value = "Passed"
$ws.Range("G24").Value = "Jonatas"
$ws.Range("E25").Value = "10 test images"
$ws.Range("F25").Value = "Passed"
$ws.Range("G25").Value = "Jonatas"
$ws.Range("E26").Value = "10 new images"
$ws.Range("F26").Value = "Passed"
$ws.Range("G26").Value = "Jon/Angelica"
$ws.Range("F27").Value = "Passed"
$ws.Range("G27").Value = "Jonatas"

# Widen the newly-populated "Assigned to" column to fit its contents
$ws.Columns.Item(7).AutoFit()
$ws.Columns.Item(7).ColumnWidth = 11.5

# Test Plan becomes the active/selected sheet & cell (was Defects before)
$ws.Activate()
$ws.Range("E25").Select()
$wsDefects.Range("F4").Select()
$ws.Activate()

Write-Host "done"
